# Replace the three-digit x one-digit multiplication problems/answers
# in the table with newly generated ones, per the commit diff.

$d = $word.ActiveDocument

$replacements = @(
    @("575×9=5175", "798×2=1596"),
    @("928×9=8352", "582×4=2328"),
    @("623×6=3738", "561×5=2805"),
    @("503×5=2515", "301×8=2408"),
    @("969×6=5814", "725×3=2175"),
    @("674×6=4044", "407×9=3663"),
    @("158×6=948",  "979×2=1958"),
    @("677×4=2708", "368×4=1472"),
    @("391×4=1564", "657×4=2628"),
    @("487×4=1948", "310×8=2480"),
    @("815×6=4890", "308×8=2464"),
    @("474×8=3792", "642×4=2568"),
    @("690×9=6210", "227×8=1816"),
    @("719×9=6471", "264×8=2112"),
    @("620×6=3720", "965×9=8685"),
    @("327×6=1962", "505×6=3030"),
    @("431×3=1293", "223×2=446"),
    @("332×9=2988", "538×9=4842"),
    @("532×4=2128", "806×2=1612"),
    @("915×7=6405", "468×9=4212"),
    @("346×4=1384", "295×4=1180"),
    @("704×2=1408", "557×3=1671"),
    @("434×9=3906", "198×5=990"),
    @("291×2=582",  "334×5=1670"),
    @("557×2=1114", "595×7=4165")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
